$d = $word.ActiveDocument
$d.Content.Find.Execute("Games:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Shipped Games:", 2)
